# Added login test case
# Fills in the first blank test-case row (row 3, under the "To Do App
# (Navigation and Login)" scenario) with a concrete login test case, and
# formats the new cells to match the rest of the filled-in rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values for the new test case row ---------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "test_<LoginwithUsernameandPassword>"
$ws.Range("C3").Value = "This is to test whether user can be logged in with the correct username and password"
$ws.Range("D3").Value = "username: ""username""`npassword: ""p@ssw0rd"""
$ws.Range("E3").Value = "Navigates to To Do page"

# --- Formatting ---------------------------------------------------------
# Row number cell: right aligned, wrapped, Calibri 10pt black
$ws.Range("A3").Font.Name = "Calibri"
$ws.Range("A3").Font.Size = 10
$ws.Range("A3").Font.Color = 0
$ws.Range("A3").HorizontalAlignment = -4152
$ws.Range("A3").WrapText = $true

# Test name cell: wrapped, Arial 10pt black
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Size = 10
$ws.Range("B3").Font.Color = 0
$ws.Range("B3").WrapText = $true

# Test description / test values cells: wrapped, Calibri 10pt black
$ws.Range("C3:D3").Font.Name = "Calibri"
$ws.Range("C3:D3").Font.Size = 10
$ws.Range("C3:D3").Font.Color = 0
$ws.Range("C3:D3").WrapText = $true

# Expected outcome cell: vertically centered, Calibri 10pt black
$ws.Range("E3").Font.Name = "Calibri"
$ws.Range("E3").Font.Size = 10
$ws.Range("E3").Font.Color = 0
$ws.Range("E3").VerticalAlignment = -4108

# Keep row 3 at its original fixed height (wrapping would otherwise grow it).
$ws.Rows(3).RowHeight = 15.75

# --- Misc ----------------------------------------------------------------
# Move the active selection, matching the author's session.
$ws.Range("F5").Select()

# Match the printer/page setup recorded when the workbook was saved.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
